$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.698.96"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "3.608.28"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.216"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.38%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000303"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "4.180.07"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "681.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +15.07%  "
$ws.Range("D16").Value = "70.758.52"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.634.71"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "110.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.885.37"
$ws.Range("E36").Value = "  +2.46%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0853"
$ws.Range("E37").Value = "  +5.88%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "514.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0468"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("B51").Value = "Jupiter"
$ws.Range("C51").Value = "https://coinranking.com/coin/qMgTxtv34+jupiter-jup"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +22.28%  "
